$d = $word.ActiveDocument

# --- Hunk 1 ---------------------------------------------------------------
# "You  are" / "proofErr gramStart..gramEnd" / " on a space station trying
# to survive" (two runs straddling a grammar-check proofErr pair) becomes a
# single plain run with the same text and no proofErr markers. The proofErr
# anchors are not addressable through the Range/text model, so we rebuild
# the paragraph: insert the merged text as a new paragraph right before the
# old one, then delete the old (still proofErr-laden) paragraph entirely.
$p2 = $d.Paragraphs(2)
$ip = $d.Range($p2.Range.Start, $p2.Range.Start)
$ip.InsertBefore("You  are on a space station trying to survive" + [char]13)
$oldP2 = $d.Paragraphs(3)
$oldP2.Range.Delete()

# --- Hunk 2 -----------------------------------------------------------------
# The last (empty, list-formatted) paragraph becomes a "Tile Sets" heading
# with a bottom border (matching the document's first paragraph style),
# followed by a blank paragraph and a paragraph containing just "s".
$last = $d.Paragraphs($d.Paragraphs.Count)

# Strip the list numbering/style and give it a bottom border instead.
$last.Range.ListFormat.RemoveNumbers()
$last.Range.ParagraphFormat.Style = "Normal"
$bottomBorder = $last.Range.ParagraphFormat.Borders.Item(-3)
$bottomBorder.LineStyle = 1
$bottomBorder.LineWidth = 2
$bottomBorder.ColorIndex = 0
$last.Range.ParagraphFormat.Borders.DistanceFromBottom = 1
$last.Range.Text = "Tile Sets"

# Append a new, blank paragraph (no border/style/run).
$endPos = $d.Content.End
$d.Range($endPos, $endPos).InsertParagraphAfter()
$blank = $d.Paragraphs($d.Paragraphs.Count)
$blank.Range.ParagraphFormat.Borders.Item(-3).LineStyle = 0
$blank.Range.Text = "X"
$blank2 = $d.Paragraphs($d.Paragraphs.Count)
$d.Range($blank2.Range.Start, $blank2.Range.Start + 1).Delete()

# Append the final paragraph containing "s".
$endPos2 = $d.Content.End
$d.Range($endPos2, $endPos2).InsertParagraphAfter()
$sPar = $d.Paragraphs($d.Paragraphs.Count)
$sPar.Range.ParagraphFormat.Borders.Item(-3).LineStyle = 0
$sPar.Range.Text = "s"
